$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.003206777572631836
$ws.Range("C2").Value = 0.01364140510559082
$ws.Range("D2").Value = 0.002009677886962891
$ws.Range("E2").Value = 0.009505367279052735
$ws.Range("F2").Value = 0.002011299133300781
$ws.Range("G2").Value = 0.008799409866333008
$ws.Range("H2").Value = 0.003192949295043945
$ws.Range("I2").Value = 0.05536923408508301
$ws.Range("J2").Value = 0.002994108200073242
$ws.Range("K2").Value = 0.01307592391967773
$ws.Range("L2").Value = 0.001986837387084961
$ws.Range("M2").Value = 0.01120638847351074
$ws.Range("B3").Value = 0.004786396026611328
$ws.Range("C3").Value = 0.009506940841674805
$ws.Range("D3").Value = 0.002593612670898437
$ws.Range("E3").Value = 0.007008838653564453
$ws.Range("F3").Value = 0.002210474014282226
$ws.Range("G3").Value = 0.006648778915405273
$ws.Range("H3").Value = 0.004607248306274414
$ws.Range("I3").Value = 0.008445501327514648
$ws.Range("J3").Value = 0.003815269470214844
$ws.Range("K3").Value = 0.008591890335083008
$ws.Range("L3").Value = 0.00279836654663086
$ws.Range("M3").Value = 0.006981229782104493
$ws.Range("B4").Value = 0.004799222946166993
$ws.Range("C4").Value = 0.008528375625610351
$ws.Range("D4").Value = 0.003587770462036133
$ws.Range("E4").Value = 0.007804679870605469
$ws.Range("F4").Value = 0.002787590026855469
$ws.Range("G4").Value = 0.007365036010742188
$ws.Range("H4").Value = 0.004198884963989258
$ws.Range("I4").Value = 0.008708524703979491
$ws.Range("J4").Value = 0.004718685150146484
$ws.Range("K4").Value = 0.008786964416503906
$ws.Range("L4").Value = 0.002201557159423828
$ws.Range("M4").Value = 0.006751203536987304
$ws.Range("B5").Value = 0.004417228698730469
$ws.Range("C5").Value = 0.008447790145874023
$ws.Range("D5").Value = 0.002785968780517578
$ws.Range("E5").Value = 0.007803583145141601
$ws.Range("H5").Value = 0.004791355133056641
$ws.Range("I5").Value = 0.009067296981811523
$ws.Range("J5").Value = 0.003607368469238281
$ws.Range("K5").Value = 0.008164691925048827
$ws.Range("B6").Value = 0.01899828910827637
$ws.Range("C6").Value = 0.01659345626831055
$ws.Range("D6").Value = 0.01160740852355957
$ws.Range("E6").Value = 0.01233224868774414
$ws.Range("F6").Value = 0.009791564941406251
$ws.Range("G6").Value = 0.01162500381469727
$ws.Range("H6").Value = 0.01860661506652832
$ws.Range("I6").Value = 0.01905021667480469
$ws.Range("J6").Value = 0.01762595176696777
$ws.Range("K6").Value = 0.01829848289489746
$ws.Range("L6").Value = 0.009198760986328125
$ws.Range("M6").Value = 0.01219363212585449
